$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.624.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.10%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.053.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.11%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'247.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.41%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +0.85%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D8").Value = "'54.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -6.92%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'60.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.93%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -2.92%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.0753"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.54%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -3.20%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.970"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +9.40%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  -4.36%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.356.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.02%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'5.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.23%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.056.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.84%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'36.554.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.18%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'17.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -5.45%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'72.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.38%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  -3.12%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'238.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.07%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'5.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.58%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +0.10%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  -2.48%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.89%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'166.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.09%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'9.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -8.18%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'20.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.05%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -1.81%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +8.33%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'5.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -7.09%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'4.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.00%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.0595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.86%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +0.00%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +1.85%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.22%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'1.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.50%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").Value = "'5.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.29%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = "'1.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.77%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -4.94%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -3.72%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -4.92%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'94.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.06%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.0915"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.80%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.416.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +8.76%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'15.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -6.15%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'7.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +10.92%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +1.59%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -3.86%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'45.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.29%  "
$ws.Range("E51").Style = "Normal"
